$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, shifting existing rows 180-223 down to 181-224.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with a new weekly price observation,
# replicating the fixed fields shared by every row in this block.
$ws.Range("A180").Value = 3
$ws.Range("B180").Value = "Femacal de La Calera"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 44543
$ws.Range("E180").Value = 5
$ws.Range("F180").Value = 100112039
$ws.Range("G180").Value = "Ciboulette"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 160
$ws.Range("K180").Value = 1500
$ws.Range("L180").Value = 1500
$ws.Range("M180").Value = 1500
$ws.Range("N180").Value = "$/docena de atados"
$ws.Range("O180").Value = "Provincia de Quillota"
$ws.Range("P180").Value = 500
$ws.Range("Q180").Value = 3
$ws.Range("R180").Value = "Hortaliza"
